# Update the Tnfsf13 -> Tnfrsf14 LR-pair sheet: "FAPs" is now included as a
# sending cluster as well as a target cluster (3x3 grid of sending/target
# clusters instead of 2x3), and all numeric stats are recomputed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @(
    @("ECs", "Tnfsf13", "Tnfrsf14", "ECs", 3, 1, 0.291812, 0.875436, 0.1179458223311005, 0.1179458223311005, 3, 1, 7.118296000000001, 21.354888, 0.4484208155031491, 0.4484208155031491, 2.077204192352, 18.694837731168, 0.05288936183490162, 0.05288936183490161),
    @("ECs", "Tnfsf13", "Tnfrsf14", "FAPs", 3, 1, 0.291812, 0.875436, 0.1179458223311005, 0.1179458223311005, 3, 1, 6.371014666666666, 19.113044, 0.4013454332903815, 0.4013454332903815, 1.859138531909333, 16.732246787184, 0.04733701716826589, 0.04733701716826588),
    @("ECs", "Tnfsf13", "Tnfrsf14", "sCs", 3, 1, 0.291812, 0.875436, 0.1179458223311005, 0.1179458223311005, 3, 1, 2.384832, 7.154496, 0.1502337512064693, 0.1502337512064693, 0.6959225955839999, 6.263303360256, 0.01771944332793299, 0.01771944332793299),
    @("FAPs", "Tnfsf13", "Tnfrsf14", "ECs", 2, 0.6666666666666666, 0.514066, 1.542198, 0.2077773947009016, 0.2077773947009016, 3, 1, 7.118296000000001, 21.354888, 0.4484208155031491, 0.4484208155031491, 3.659273951536, 32.933465563824, 0.093171708774898, 0.09317170877489798),
    @("FAPs", "Tnfsf13", "Tnfrsf14", "FAPs", 2, 0.6666666666666666, 0.514066, 1.542198, 0.2077773947009016, 0.2077773947009016, 3, 1, 6.371014666666666, 19.113044, 0.4013454332903815, 0.4013454332903815, 3.275122025634667, 29.476098230712, 0.08339050850417998, 0.08339050850417998),
    @("FAPs", "Tnfsf13", "Tnfrsf14", "sCs", 2, 0.6666666666666666, 0.514066, 1.542198, 0.2077773947009016, 0.2077773947009016, 3, 1, 2.384832, 7.154496, 0.1502337512064693, 0.1502337512064693, 1.225961046912, 11.033649422208, 0.03121517742182364, 0.03121517742182364),
    @("sCs", "Tnfsf13", "Tnfrsf14", "ECs", 3, 1, 1.668241, 5.004723, 0.674276782967998, 0.674276782967998, 3, 1, 7.118296000000001, 21.354888, 0.4484208155031491, 0.4484208155031491, 11.875033237336, 106.875299136024, 0.3023597448933495, 0.3023597448933495),
    @("sCs", "Tnfsf13", "Tnfrsf14", "FAPs", 3, 1, 1.668241, 5.004723, 0.674276782967998, 0.674276782967998, 3, 1, 6.371014666666666, 19.113044, 0.4013454332903815, 0.4013454332903815, 10.62838787853467, 95.65549090681199, 0.2706179076179357, 0.2706179076179357),
    @("sCs", "Tnfsf13", "Tnfrsf14", "sCs", 3, 1, 1.668241, 5.004723, 0.674276782967998, 0.674276782967998, 3, 1, 2.384832, 7.154496, 0.1502337512064693, 0.1502337512064693, 3.978474520512, 35.806270684608, 0.1012991304567127, 0.1012991304567127)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowValues = $rows[$i]
    for ($c = 0; $c -lt $cols.Count; $c++) {
        $ws.Range($cols[$c] + $r).Value = $rowValues[$c]
    }
}
